$wb = $excel.ActiveWorkbook

# Insert the new "table_definitions" worksheet before the current first sheet
# (Excel renumbers sheetId/r:id for everything after it, exactly as in the diff).
$firstSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($firstSheet)
$ws.Name = "table_definitions"

# --- Header row ---
$ws.Range("A1").Value = "mapping_file_name"
$ws.Range("B1").Value = "entity_name"
$ws.Range("C1").Value = "required_entities"
$ws.Range("D1").Value = "destination_table_name"
$ws.Range("E1").Value = "table_type"
$ws.Range("F1").Value = "source_table_name"
$ws.Range("G1").Value = "casrec_conditions"
$ws.Range("H1").Value = "source_table_additional_columns"

# --- Data rows ---
$ws.Range("A2").Value = "client_nodebtchase_warnings"
$ws.Range("B2").Value = "warnings"
$ws.Range("D2").Value = "warnings"
$ws.Range("E2").Value = "data"
$ws.Range("F2").Value = "pat"
$ws.Range("G2").Value = "Debt chase = not null"
$ws.Range("H2").Value = "Case"

$ws.Range("A3").Value = "client_saarcheck_warnings"
$ws.Range("B3").Value = "warnings"
$ws.Range("D3").Value = "warnings"
$ws.Range("E3").Value = "data"
$ws.Range("F3").Value = "pat"
$ws.Range("G3").Value = "SAAR Check = not null"
$ws.Range("H3").Value = "Case"

$ws.Range("A4").Value = "client_special_warnings"
$ws.Range("B4").Value = "warnings"
$ws.Range("D4").Value = "warnings"
$ws.Range("E4").Value = "data"
$ws.Range("F4").Value = "pat"
$ws.Range("G4").Value = "SIM = not null"
$ws.Range("H4").Value = "Case"

$ws.Range("A5").Value = "client_violent_warnings"
$ws.Range("B5").Value = "warnings"
$ws.Range("D5").Value = "warnings"
$ws.Range("E5").Value = "data"
$ws.Range("F5").Value = "pat"
$ws.Range("G5").Value = "VWM = not null"
$ws.Range("H5").Value = "Case"

$ws.Range("A6").Value = "client_person_warning"
$ws.Range("B6").Value = "warnings"
$ws.Range("D6").Value = "warnings"
$ws.Range("E6").Value = "join"

$ws.Range("A7").Value = "deputy_special_warnings"
$ws.Range("B7").Value = "warnings"
$ws.Range("D7").Value = "warnings"
$ws.Range("E7").Value = "data"
$ws.Range("F7").Value = "deputy"
$ws.Range("G7").Value = "SIM = not null"
$ws.Range("H7").Value = "Deputy No"

$ws.Range("A8").Value = "deputy_violent_warnings"
$ws.Range("B8").Value = "warnings"
$ws.Range("D8").Value = "warnings"
$ws.Range("E8").Value = "data"
$ws.Range("F8").Value = "deputy"
$ws.Range("G8").Value = "VWM = not null"
$ws.Range("H8").Value = "Deputy No"

$ws.Range("A9").Value = "deputy_person_warning"
$ws.Range("B9").Value = "warnings"
$ws.Range("D9").Value = "warnings"
$ws.Range("E9").Value = "join"

# --- Formatting: header row uses a grey Helvetica font ---
# (Data rows render in the workbook's default font - Arial 10pt black - which is
# already the implicit style, so no explicit font assignment is needed there.)
$headerRange = $ws.Range("A1:H1")
$headerRange.Font.Name = "Helvetica"
$headerRange.Font.Size = 10
$headerRange.Font.Color = 3355443

# Make the new sheet the active / selected tab, with the whole used range selected
$ws.Activate()
$ws.Range("A1:H9").Select()
